$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Section "1.1." -> "1."  (top of document header block)
# ------------------------------------------------------------------
$d.Content.Find.Execute("1.1.", $false, $false, $false, $false, $false, $true, 1, $false, "1.", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "Thận trọng trong việc bảo quản" -> "Thận trọng trong xử lý an toàn"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Thận trọng trong việc bảo quản", $false, $false, $false, $false, $false, $true, 1, $false, "Thận trọng trong xử lý an toàn", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Remove the bookmark around "Thường xuyên kiểm tra rò rỉ." and
#    merge the preceding space into the same run as the sentence.
# ------------------------------------------------------------------
$d.Content.Find.Execute(" Thường xuyên kiểm tra rò rỉ.", $false, $false, $false, $false, $false, $true, 1, $false, " Thường xuyên kiểm tra rò rỉ.", 2) | Out-Null
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# 4) "Hơi tương đối" -> "Tỉ trọng hơi"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Hơi tương đối", $false, $false, $false, $false, $false, $true, 1, $false, "Tỉ trọng hơi", 2) | Out-Null

# ------------------------------------------------------------------
# 5) "Tính chất oxy hóa" -> "Tính oxy hóa"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Tính chất oxy hóa", $false, $false, $false, $false, $false, $true, 1, $false, "Tính oxy hóa", 2) | Out-Null

# ------------------------------------------------------------------
# 6) "Cacbon dioxide" -> "Cacbon đioxit" (written lowercase "bon dioxide")
# ------------------------------------------------------------------
$d.Content.Find.Execute("bon dioxide", $false, $false, $false, $false, $false, $true, 1, $false, "bon đioxit", 2) | Out-Null

# ------------------------------------------------------------------
# 7) "cacbon monoxide" -> "cacbon monoxit"
# ------------------------------------------------------------------
$d.Content.Find.Execute("bon monoxide", $false, $false, $false, $false, $false, $true, 1, $false, "bon monoxit", 2) | Out-Null

# ------------------------------------------------------------------
# 8) "Độc tính cấp (đường miệng)" -> "Độc cấp tính (đường miệng)"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Độc tính cấp (đường miệng)", $false, $false, $false, $false, $false, $true, 1, $false, "Độc cấp tính (đường miệng)", 2) | Out-Null

# ------------------------------------------------------------------
# 9) "Độc tính cấp (qua da)" -> "Độc cấp tính (qua da)"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Độc tính cấp (qua da)", $false, $false, $false, $false, $false, $true, 1, $false, "Độc cấp tính (qua da)", 2) | Out-Null

# ------------------------------------------------------------------
# 10) "Độc tính cấp tính đối với thủy sinh: " -> "Độc tính cấp tính đối với sinh vật thủy sinh: "
# ------------------------------------------------------------------
$d.Content.Find.Execute("Độc tính cấp tính đối với thủy sinh: ", $false, $false, $false, $false, $false, $true, 1, $false, "Độc tính cấp tính đối với sinh vật thủy sinh: ", 2) | Out-Null

# ------------------------------------------------------------------
# 11) "Độc tính mãn tính đối với thủy sinh:" -> "Độc tính mãn tính đối với sinh vật thủy sinh:"
#     plus a fresh "_GoBack" bookmark placed right after "sinh vật "
# ------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Độc tính mãn tính đối với thủy sinh:")
if ($found) {
    $prefix = "Độc tính mãn tính đối với "
    $insStart = $r.Start + $prefix.Length
    $ins = $d.Range($insStart, $insStart)
    $ins.InsertBefore("sinh vật ")
    $bmPos = $insStart + ("sinh vật ").Length
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
